$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.519.73"
$ws.Range("E2").Value = "  +5.87%  "

$ws.Range("D3").Value = "3.572.35"
$ws.Range("E3").Value = "  +5.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.51"
$ws.Range("E5").Value = "  +5.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.94"
$ws.Range("E6").Value = "  +8.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +1.65%  "

$ws.Range("D8").Value = "3.566.12"
$ws.Range("E8").Value = "  +5.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  +5.97%  "

$ws.Range("E11").Value = "  +3.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.24"
$ws.Range("E12").Value = "  +8.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000292"
$ws.Range("E13").Value = "  +5.07%  "

$ws.Range("E14").Value = "  +5.16%  "

$ws.Range("D15").Value = "4.135.63"
$ws.Range("E15").Value = "  +5.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.27"
$ws.Range("E16").Value = "  +5.24%  "

$ws.Range("D17").Value = "3.561.80"
$ws.Range("E17").Value = "  +4.56%  "

$ws.Range("D18").Value = "69.457.31"
$ws.Range("E18").Value = "  +5.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.43"
$ws.Range("E19").Value = "  +4.78%  "

$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("E21").Value = "  +4.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "503.15"
$ws.Range("E22").Value = "  +4.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.54"
$ws.Range("E23").Value = "  +11.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.12"
$ws.Range("E24").Value = "  +19.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.44"
$ws.Range("E25").Value = "  +8.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.19"
$ws.Range("E26").Value = "  +1.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.06"
$ws.Range("E27").Value = "  +4.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.16"
$ws.Range("E28").Value = "  +4.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  +6.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.05"
$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.51"
$ws.Range("E31").Value = "  +14.42%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "618.99"
$ws.Range("E32").Value = "  +7.28%  "

$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.15"
$ws.Range("E33").Value = "  +5.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.32"
$ws.Range("E34").Value = "  +4.02%  "

$ws.Range("E35").Value = "  +6.42%  "

$ws.Range("D36").Value = "0.0₃0835"
$ws.Range("E36").Value = "  +12.68%  "

$ws.Range("E37").Value = "  +4.76%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.97"
$ws.Range("E39").Value = "  +5.80%  "

$ws.Range("E40").Value = "  +6.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").Value = "3.331.32"
$ws.Range("E42").Value = "  +7.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.06"
$ws.Range("E43").Value = "  +9.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.72"
$ws.Range("E44").Value = "  +11.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0442"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  +21.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  +3.87%  "

$ws.Range("E48").Value = "  +2.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("E49").Value = "  +7.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("E50").Value = "  +4.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.17%  "
